$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($CellAddress, $Value) {
    $r = $ws.Range($CellAddress)
    $r.NumberFormat = "@"
    $r.Value = $Value
    $r.Style = "Normal"
}

Set-TextValue 'D2' '59.237.85'
$ws.Range('E2').Value = '  +0.86%  '
Set-TextValue 'D3' '2.585.21'
$ws.Range('E3').Value = '  -0.08%  '
Set-TextValue 'D4' '1.00'
$ws.Range('E4').Value = '  +0.05%  '
Set-TextValue 'D5' '567.07'
$ws.Range('E5').Value = '  +2.42%  '
Set-TextValue 'D6' '143.59'
$ws.Range('E6').Value = '  +0.15%  '
Set-TextValue 'D7' '0.999'
$ws.Range('E7').Value = '  +0.06%  '
Set-TextValue 'D8' '0.600'
$ws.Range('E8').Value = '  -0.22%  '
Set-TextValue 'D9' '2.595.06'
$ws.Range('E9').Value = '  +0.06%  '
Set-TextValue 'D10' '6.66'
$ws.Range('E10').Value = '  -1.61%  '
$ws.Range('E11').Value = '  +3.19%  '
Set-TextValue 'D12' '0.154'
$ws.Range('E12').Value = '  +9.58%  '
Set-TextValue 'D13' '0.343'
$ws.Range('E13').Value = '  +2.32%  '
Set-TextValue 'D14' '3.043.49'
$ws.Range('E14').Value = '  -0.03%  '
Set-TextValue 'D15' '59.286.37'
$ws.Range('E15').Value = '  +1.10%  '
Set-TextValue 'D16' '22.42'
$ws.Range('E16').Value = '  +8.08%  '
$ws.Range('E17').Value = '  +4.47%  '
Set-TextValue 'D18' '2.593.36'
$ws.Range('E18').Value = '  -0.06%  '
Set-TextValue 'D19' '4.52'
$ws.Range('E19').Value = '  +1.91%  '
Set-TextValue 'D20' '336.23'
$ws.Range('E20').Value = '  +0.02%  '
Set-TextValue 'D21' '10.21'
$ws.Range('E21').Value = '  +1.91%  '
$ws.Range('E22').Value = '  +1.43%  '
Set-TextValue 'D24' '64.20'
$ws.Range('E24').Value = '  -3.02%  '
Set-TextValue 'D25' '0.450'
$ws.Range('E25').Value = '  +6.09%  '
$ws.Range('E26').Value = '  -0.16%  '
$ws.Range('E27').Value = '  +2.38%  '
Set-TextValue 'D28' '7.25'
$ws.Range('E28').Value = '  +2.16%  '
Set-TextValue 'D29' '0.0₃0785'
$ws.Range('E29').Value = '  +4.33%  '
$ws.Range('E30').Value = '  -0.07%  '
$ws.Range('E31').Value = '  +0.40%  '
$ws.Range('B32').Value = 'Aptos'
$ws.Range('C32').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue 'D32' '6.04'
$ws.Range('E32').Value = '  +2.16%  '
$ws.Range('B33').Value = 'Monero'
$ws.Range('C33').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue 'D33' '157.81'
$ws.Range('E33').Value = '  +2.53%  '
Set-TextValue 'D34' '19.02'
$ws.Range('E34').Value = '  +0.88%  '
Set-TextValue 'D35' '4.04'
$ws.Range('E35').Value = '  +3.45%  '
Set-TextValue 'D36' '0.887'
$ws.Range('B37').Value = 'SuiNetwork'
$ws.Range('C37').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
Set-TextValue 'D37' '0.878'
$ws.Range('E37').Value = '  +1.60%  '
$ws.Range('B38').Value = 'ImmutableX'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue 'D38' '1.14'
$ws.Range('E38').Value = '  +2.76%  '
$ws.Range('E39').Value = '  +3.88%  '
Set-TextValue 'D40' '36.91'
$ws.Range('E40').Value = '  -0.09%  '
$ws.Range('B41').Value = 'Bittensor'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue 'D41' '294.25'
$ws.Range('E41').Value = '  +4.05%  '
$ws.Range('B42').Value = 'Filecoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 'D42' '3.65'
$ws.Range('E42').Value = '  +1.60%  '
Set-TextValue 'D43' '0.998'
$ws.Range('E43').Value = '  +0.09%  '
Set-TextValue 'D44' '0.0975'
$ws.Range('E44').Value = '  +1.91%  '
$ws.Range('E45').Value = '  +0.35%  '
$ws.Range('E46').Value = '  +1.16%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D47' '19.21'
$ws.Range('E47').Value = '  +3.26%  '
$ws.Range('B48').Value = 'WhiteBITCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
Set-TextValue 'D48' '10.62'
$ws.Range('E48').Value = '  -0.13%  '
$ws.Range('B49').Value = 'VeChain'
$ws.Range('C49').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D49' '0.0233'
$ws.Range('E49').Value = '  +2.90%  '
$ws.Range('B50').Value = 'Aave'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue 'D50' '124.29'
$ws.Range('E50').Value = '  +7.82%  '
Set-TextValue 'D51' '18.56'
$ws.Range('E51').Value = '  +4.90%  '
